$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") from row 2 to row 101: 45243 -> 45244
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
